$d = $word.ActiveDocument

$replacements = @(
    @("83÷6=", "83÷2="),
    @("49÷8=", "35÷2="),
    @("43÷5=", "13÷6="),
    @("93÷6=", "92÷6="),
    @("50÷6=", "44÷8="),
    @("34÷9=", "17÷8="),
    @("88÷4=", "15÷3="),
    @("86÷7=", "27÷7="),
    @("76÷9=", "49÷4="),
    @("29÷4=", "17÷3="),
    @("90÷9=", "35÷4="),
    @("97÷2=", "95÷3="),
    @("84÷7=", "57÷3="),
    @("53÷3=", "30÷9="),
    @("78÷6=", "87÷6="),
    @("80÷9=", "45÷8="),
    @("37÷6=", "71÷6="),
    @("13÷5=", "22÷8="),
    @("82÷5=", "65÷6="),
    @("82÷6=", "33÷9="),
    @("74÷6=", "59÷3="),
    @("82÷9=", "90÷5="),
    @("81÷7=", "42÷2="),
    @("91÷9=", "25÷9="),
    @("71÷5=", "19÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
